$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.869.98"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "1.917.91"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.99"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4565"
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3806"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07753"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.25"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("D12").Value = "1.938.62"
$ws.Range("E12").Value = "  +3.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.993"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06972"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.38"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009499"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.64"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "28.876.05"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.344"
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.11"
$ws.Range("E23").Value = "  +2.29%  "
$ws.Range("D24").Value = "2.155.32"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.067"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.04"
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.08"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.623"
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.96"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.853"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09277"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8708"
$ws.Range("E32").Value = "  +1.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.101"
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.246"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.026"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05709"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.148"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.002"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02040"
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.064"
$ws.Range("E40").Value = "  +11.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.506"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5506"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1759"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.356"
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000002874"
$ws.Range("E45").Value = "  +16.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.172"
$ws.Range("E46").Value = "  +3.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5161"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06925"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.11"
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.63"
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.765"
$ws.Range("E51").Value = "  -0.57%  "
